$p = $ppt.ActivePresentation

# --- Slide 7: "Second round match percentage" textbox (TextBox 7 / shape 6) ---
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(6)
$tr7 = $shp7.TextFrame.TextRange
$tr7.Characters(1, $tr7.Length).Text = "Second round match percentage = 20% for WB, 25% for TRF"

# --- Slide 8: "Third round match percentage" textbox (TextBox 7 / shape 6) ---
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(6)
$tr8 = $shp8.TextFrame.TextRange
$tr8.Characters(1, $tr8.Length).Text = "Third round match percentage = 4% for WB, 9% for TRF"

# --- Slide 9: "Fourth round match percentage = 1% ..." textbox (TextBox 7 / shape 6) ---
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(6)
$tr9 = $shp9.TextFrame.TextRange
$tr9.Characters(1, $tr9.Length).Text = "Fourth round match percentage = 1% for WB, 2% for TRF "

# --- Slide 10: "Fourth round match percentage = 50% ..." textbox (TextBox 7 / shape 6) ---
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(6)
$tr10 = $shp10.TextFrame.TextRange
$tr10.Characters(1, $tr10.Length).Text = "Fourth round match percentage = 50% for WB, 25% for TRF"
